$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.588.85"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.220.30"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.64"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.08"
$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.20"
$ws.Range("E10").Value = "  -1.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -2.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.87"
$ws.Range("E12").Value = "  -3.13%  "

$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.555.96"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.67"
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.222.98"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("E18").Value = "  -3.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.449.39"
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.69"
$ws.Range("E21").Value = "  -2.08%  "

$ws.Range("E22").Value = "  -4.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.00"
$ws.Range("E23").Value = "  -9.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "229.02"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  +6.26%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("E27").Value = "  -3.75%  "

$ws.Range("E28").Value = "  -5.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -2.41%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.64"
$ws.Range("E31").Value = "  +3.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.14"
$ws.Range("E32").Value = "  +17.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.28"
$ws.Range("E33").Value = "  -1.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  -3.56%  "

$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.38"
$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("E39").Value = "  +6.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.46"
$ws.Range("E40").Value = "  -4.16%  "

$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.49"
$ws.Range("E42").Value = "  -2.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.62"
$ws.Range("E43").Value = "  -5.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.197"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0987"
$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.28"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.423"
$ws.Range("E51").Value = "  +14.65%  "
